$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.677.12"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "3.784.13"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.75"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.08"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").Value = "3.781.64"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  -0.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.450"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.73"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "4.422.25"
$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").Value = "3.774.49"
$ws.Range("E16").Value = "  -0.66%  "

$ws.Range("D17").Value = "68.661.35"
$ws.Range("E17").Value = "  +1.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.37"
$ws.Range("E18").Value = "  -3.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.09"
$ws.Range("E19").Value = "  -2.52%  "

$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.44"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.704"
$ws.Range("E23").Value = "  -3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.72"
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("E25").Value = "  -2.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.22"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.32"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "3.931.73"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("E31").Value = "  -4.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.50"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.27"
$ws.Range("E34").Value = "  -0.83%  "

$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"

$ws.Range("D37").Value = "3.739.64"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("E39").Value = "  -7.35%  "

$ws.Range("E40").Value = "  +0.55%  "

$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.85"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.309"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.68"
$ws.Range("E47").Value = "  +12.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.61"
$ws.Range("E48").Value = "  -1.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "407.06"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.52"
$ws.Range("E50").Value = "  -1.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.91"
$ws.Range("E51").Value = "  +2.07%  "
